# Update the phones worksheet:
#  - insert two new leading columns ("president", "dog") before the
#    existing "client"/"calls" columns (renamed "pclient"/"pcalls")
#  - populate the new columns with boolean-style (0/1) sample data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the original "client"/"calls" data (columns A:B) two columns to
# the right so it lands in C:D, preserving the original values/order.
$ws.Range("A1:B4").Cut($ws.Range("C1:D4"))

# Rename the shifted header cells.
$ws.Range("C1").Value = "pclient"
$ws.Range("D1").Value = "pcalls"

# New leading columns with header + sample 0/1 data.
$ws.Range("A1").Value = "president"
$ws.Range("B1").Value = "dog"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 0

$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 1
